$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed entirely ("RM 232" originally at row 26,
# "SC 92" originally at row 28). Delete the lower row first so the row index of
# the upper one doesn't shift before we get to it.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Cell-level value updates (imputed / cleared values), using the row numbers
# that result AFTER the two rows above have been removed.
$ws.Range("C2").Value = 14.9
$ws.Range("C3").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("C11").Value = 11.4
$ws.Range("C13").Value = ""
$ws.Range("C21").Value = 12.7
$ws.Range("C25").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
